$d = $word.ActiveDocument

# ============================================================
# Change 1: insert a new "Meta description" paragraph right
# after the title paragraph ("Play Asgardian Stones Slot Game
# for Free - A NetEnt Journey to Norse Mythology"), and before
# the "Asgardian Stones Slot Game: Journey Into Norse
# Mythology" Heading2 paragraph.
# ============================================================

# Split the document right before paragraph 2 (the Heading2
# paragraph) so a brand-new, empty paragraph is created in
# between. This new paragraph is given the leading empty run
# that every other body paragraph in this document has.
$titleHeadingPara = $d.Paragraphs(2)
$splitPoint = $d.Range($titleHeadingPara.Range.Start, $titleHeadingPara.Range.Start)
$splitPoint.InsertParagraphBefore()

# Grab the formatted text of the (soon to be removed) bold
# "Play Asgardian Stones..." paragraph near the end of the
# document - it already has the exact run shape we need
# (an empty run followed by one bold run), so reusing it lets
# us reproduce that structure faithfully instead of typing
# fresh text into the new paragraph.
$paraCount = $d.Paragraphs.Count
$trailingBoldPara = $d.Paragraphs($paraCount - 1)
$boldFormattedText = $trailingBoldPara.Range.FormattedText

$newPara = $d.Paragraphs(2)
$newPara.Style = "Normal"
$newPara.Range.FormattedText = $boldFormattedText

# Turn the copied "Play Asgardian Stones..." text into the
# bold "Meta description" label.
$newPara = $d.Paragraphs(2)
$newPara.Range.Find.Execute(
    "Play Asgardian Stones Slot Game for Free - A NetEnt Journey to Norse Mythology",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Meta description", 2)

# Append the (non-bold) rest of the meta description after the
# bold label, as a separate run.
$newPara = $d.Paragraphs(2)
$insertStart = $newPara.Range.End - 1
$insertPoint = $d.Range($insertStart, $insertStart)
$insertPoint.InsertAfter(": Read our review of Asgardian Stones Slot Game. Play for free and journey into Norse mythology with cascading symbols, bonus wheel and high-paying wins.")
$insertEnd = $d.Paragraphs(2).Range.End - 1
$newRunRange = $d.Range($insertStart, $insertEnd)
$newRunRange.Font.Bold = 0

# ============================================================
# Change 2: remove the duplicated bold title paragraph near the
# end of the document, and replace the text of the italic
# "Read our review..." paragraph (the old meta description)
# with the new feature-image prompt text.
# ============================================================

$paraCount = $d.Paragraphs.Count
$trailingBoldPara = $d.Paragraphs($paraCount - 1)
$trailingBoldPara.Range.Delete()

$imagePromptPara = $d.Paragraphs($d.Paragraphs.Count)
$imagePromptPara.Range.Find.Execute(
    "Read our review of Asgardian Stones Slot Game. Play for free and journey into Norse mythology with cascading symbols, bonus wheel and high-paying wins.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Create an eye-catching feature image for Asgardian Stones that features a happy Maya warrior with glasses. The image should be in cartoon style and should convey the excitement and thrill of the game. You could include elements from Norse mythology, such as the Asgardian Stones or symbols of power and strength to add to the theme. Make the image bright and colorful with bold outlines to make it pop. The Maya warrior should be depicted as having fun and enjoying the game to encourage players to give it a try. The image should be of high quality and clearly convey the message that this is a fun and exciting game to play.",
    2)

Write-Output "Edit complete"
